$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (pt_max) values from 50 to 60 for data rows 2 through 12
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 6).Value = 60
}
